$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh per-coin Price (column D) and Volume/1h change (column E) figures.
# Rows 31/32 also swap Coin name + Link + Price + Volume (ranking reordered).
# Numeric-looking Price strings are forced back to text so they keep their
# original formatting (e.g. "8.00", "1.90") instead of being coerced to numbers,
# then the style is reset to Normal so no visual/style change is introduced.

$ws.Range("D2").Value = "62.234.99"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.447.83"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "2.443.18"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "2.868.34"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "62.108.44"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "2.436.83"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E24").Value = "  -5.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "602.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.57%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.136"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "0.0₆0269"
$ws.Range("E48").Value = "  +21.08%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
